# Updates cryptos list values (prices / 1h volume %, and two name/link row swaps)
# to match the refreshed data snapshot from the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "57.428.85"
$ws.Range("E2").Value = "  -0.57%  "
# Row 3
$ws.Range("D3").Value = "3.100.02"
$ws.Range("E3").Value = "  +1.29%  "
# Row 4
$ws.Range("E4").Value = "  -0.02%  "
# Row 5
$ws.Range("D5").Value = "'522.64"
$ws.Range("E5").Value = "  +0.97%  "
# Row 6
$ws.Range("D6").Value = "'140.85"
$ws.Range("E6").Value = "  -0.78%  "
# Row 7
$ws.Range("E7").Value = "  -0.02%  "
# Row 8
$ws.Range("D8").Value = "3.101.17"
$ws.Range("E8").Value = "  +1.47%  "
# Row 9
$ws.Range("D9").Value = "'0.435"
$ws.Range("E9").Value = "  +0.14%  "
# Row 11
$ws.Range("E11").Value = "  +0.74%  "
# Row 12
$ws.Range("D12").Value = "'0.384"
$ws.Range("E12").Value = "  +1.78%  "
# Row 13
$ws.Range("D13").Value = "3.634.29"
$ws.Range("E13").Value = "  +1.34%  "
# Row 14
$ws.Range("D14").Value = "'0.131"
$ws.Range("E14").Value = "  +1.32%  "
# Row 15
$ws.Range("D15").Value = "'26.10"
$ws.Range("E15").Value = "  +0.63%  "
# Row 16
$ws.Range("E16").Value = "  +0.60%  "
# Row 17
$ws.Range("D17").Value = "57.536.56"
$ws.Range("E17").Value = "  -0.40%  "
# Row 18
$ws.Range("D18").Value = "3.101.05"
$ws.Range("E18").Value = "  +1.66%  "
# Row 19
$ws.Range("E19").Value = "  +0.21%  "
# Row 20
$ws.Range("D20").Value = "'12.78"
$ws.Range("E20").Value = "  -0.14%  "
# Row 21
$ws.Range("E21").Value = "  -0.63%  "
# Row 22
$ws.Range("D22").Value = "'336.09"
$ws.Range("E22").Value = "  +1.77%  "
# Row 23
$ws.Range("E23").Value = "  +0.06%  "
# Row 24
$ws.Range("D24").Value = "'0.512"
$ws.Range("E24").Value = "  +2.74%  "
# Row 25
$ws.Range("D25").Value = "'66.50"
$ws.Range("E25").Value = "  +1.21%  "
# Row 26
$ws.Range("D26").Value = "'0.168"
$ws.Range("E26").Value = "  -0.65%  "
# Row 27
$ws.Range("E27").Value = "  +0.30%  "
# Row 28
$ws.Range("D28").Value = "0.0₃0918"
$ws.Range("E28").Value = "  +1.89%  "
# Row 29
$ws.Range("D29").Value = "'6.49"
$ws.Range("E29").Value = "  +1.70%  "
# Row 30
$ws.Range("E30").Value = "  +0.01%  "
# Row 31
$ws.Range("D31").Value = "'7.20"
$ws.Range("E31").Value = "  -0.08%  "
# Row 32
$ws.Range("E32").Value = "  +2.11%  "
# Row 33
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").Value = "'20.88"
$ws.Range("E33").Value = "  +1.09%  "
# Row 34
$ws.Range("B34").Value = "Fetch.AI"
$ws.Range("C34").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D34").Value = "'1.20"
$ws.Range("E34").Value = "  +0.47%  "
# Row 35
$ws.Range("D35").Value = "'156.71"
$ws.Range("E35").Value = "  +1.32%  "
# Row 36
$ws.Range("D36").Value = "'4.64"
$ws.Range("E36").Value = "  +3.23%  "
# Row 37
$ws.Range("D37").Value = "'6.10"
$ws.Range("E37").Value = "  +2.89%  "
# Row 38
$ws.Range("D38").Value = "'26.93"
$ws.Range("E38").Value = "  +0.06%  "
# Row 39
$ws.Range("E39").Value = "  +1.53%  "
# Row 40
$ws.Range("D40").Value = "'0.0659"
$ws.Range("E40").Value = "  -2.09%  "
# Row 41
$ws.Range("B41").Value = "RenzoRestakedETH"
$ws.Range("C41").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D41").Value = "3.141.49"
$ws.Range("E41").Value = "  +1.23%  "
# Row 42
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").Value = "'3.93"
$ws.Range("E42").Value = "  +0.56%  "
# Row 43
$ws.Range("E43").Value = "  +4.71%  "
# Row 44
$ws.Range("D44").Value = "'1.51"
$ws.Range("E44").Value = "  +11.30%  "
# Row 45
$ws.Range("E45").Value = "  +0.77%  "
# Row 47
$ws.Range("D47").Value = "2.298.78"
$ws.Range("E47").Value = "  +1.90%  "
# Row 48
$ws.Range("E48").Value = "  +0.33%  "
# Row 49
$ws.Range("E49").Value = "  +5.25%  "
# Row 50
$ws.Range("D50").Value = "'20.71"
$ws.Range("E50").Value = "  +0.36%  "
# Row 51
$ws.Range("E51").Value = "  +2.31%  "
